# Extend the indicator table with two more year columns: 2021 (M) and 2022 (N),
# mirroring the existing 2020 column (L) for both the header row (3) and the
# data row (4). Then move the active selection to N15, matching the saved
# workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone column L's formatting (header style + data style) into the two new
# columns before filling in their own values.
$ws.Range("L3:L4").Copy($ws.Range("M3:M4")) | Out-Null
$ws.Range("L3:L4").Copy($ws.Range("N3:N4")) | Out-Null

# Header row: years 2021 and 2022
$ws.Range("M3").Value = 2021
$ws.Range("N3").Value = 2022

# Data row: same indicator value as the 2020 column (6.18)
$ws.Range("M4").Value = 6.18
$ws.Range("N4").Value = 6.18

# Match the saved selection state
$ws.Range("N15").Select() | Out-Null
